$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Row 53 (write new strings in the order they appear in the target sharedStrings.xml)
$ws.Range("C53").Value = "Min Max deductible over-under limit scenario"
$ws.Range("E53").Value = "1,12,10,8"
$ws.Range("B53").Value = "fm49"
$ws.Range("D53").Value = "All"
$ws.Range("F53").Value = 3
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "complete"
$ws.Range("I53").Value = "complete"

# Row 54
$ws.Range("B54").Value = "fm50"
$ws.Range("C54").Value = "OED spec example 4 - nested sublimits "
$ws.Range("D54").Value = "All"
$ws.Range("H54").Value = "in progress"
$ws.Range("I54").Value = "in progress"

# Apply styles matching row 52 (B,C,H,I -> s=5 font; G -> s=7 right aligned)
$ws.Range("B52").Copy() | Out-Null
$ws.Range("B53").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C52").Copy() | Out-Null
$ws.Range("C53").PasteSpecial(-4122) | Out-Null
$ws.Range("G52").Copy() | Out-Null
$ws.Range("G53").PasteSpecial(-4122) | Out-Null
$ws.Range("H52").Copy() | Out-Null
$ws.Range("H53").PasteSpecial(-4122) | Out-Null
$ws.Range("I52").Copy() | Out-Null
$ws.Range("I53").PasteSpecial(-4122) | Out-Null

$ws.Range("B52").Copy() | Out-Null
$ws.Range("B54").PasteSpecial(-4122) | Out-Null
$ws.Range("C52").Copy() | Out-Null
$ws.Range("C54").PasteSpecial(-4122) | Out-Null
$ws.Range("H52").Copy() | Out-Null
$ws.Range("H54").PasteSpecial(-4122) | Out-Null
$ws.Range("I52").Copy() | Out-Null
$ws.Range("I54").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("I54").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
